$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")

# Row 62
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 5235.6665
$arr[0,1] = 4149.1113
$arr[0,2] = 6322.222
$arr[0,3] = 4149.1113
$arr[0,4] = 6322.222
$arr[0,5] = -3525.1113
$arr[0,6] = -7570.222
$ws.Range("H62:N62").Value = $arr

# Row 65
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 5235.6665
$arr[0,1] = 4149.1113
$arr[0,2] = 6322.222
$arr[0,3] = 20745.5565
$arr[0,4] = 31611.11
$arr[0,5] = -17625.5565
$arr[0,6] = -37851.11
$ws.Range("H65:N65").Value = $arr

# Row 86
$ws.Range("H86").Value = 5825.269
$ws.Range("I86").Value = 3807.1875
$ws.Range("K86").Value = 3807.1875
$ws.Range("M86").Value = -2684.1875

# Row 89
$ws.Range("H89").Value = 5825.269
$ws.Range("I89").Value = 3807.1875
$ws.Range("K89").Value = 19035.9375
$ws.Range("M89").Value = -13419.9375

# Row 92
$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 683.8
$arr[0,1] = 676.8570999999999
$arr[0,2] = 700
$arr[0,3] = 676.8570999999999
$arr[0,4] = 700
$arr[0,5] = 571.1429000000001
$ws.Range("H92:M92").Value = $arr
$ws.Range("N92").Value = -3196

# Row 111
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1246.2
$arr[0,1] = 822.875
$arr[0,2] = 1730
$arr[0,3] = 2468.625
$arr[0,4] = 5190
$arr[0,5] = 598.375
$arr[0,6] = -11324
$ws.Range("H111:N111").Value = $arr

# Row 112
$ws.Range("H112").Value = 435846.44
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 435846.44
$ws.Range("K112").Value = 0
$ws.Range("N112").Value = -1309755.32
$ws.Range("M112").ClearContents()

# Row 118
$ws.Range("H118").Value = 948.9375
$ws.Range("I118").Value = 948.9375
$ws.Range("K118").Value = 2846.8125
$ws.Range("M118").Value = -1189.8125

# Row 127
$ws.Range("H127").Value = 806
$ws.Range("I127").Value = 718.55554
$ws.Range("K127").Value = 2155.66662
$ws.Range("M127").Value = 2804.33338

# Row 132
$ws.Range("H132").Value = 1159.9512
$ws.Range("I132").Value = 1176.7179
$ws.Range("K132").Value = 3530.1537
$ws.Range("M132").Value = -1000.1537

# Row 137
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 422379.25
$arr[0,1] = 2106.2942
$arr[0,2] = 842652.25
$arr[0,3] = 6318.882599999999
$arr[0,4] = 2527956.75
$arr[0,5] = -3768.882599999999
$arr[0,6] = -2533056.75
$ws.Range("H137:N137").Value = $arr

# Row 138
$ws.Range("H138").Value = 2905.459
$ws.Range("J138").Value = 3346.8462
$ws.Range("L138").Value = 10040.5386
$ws.Range("N138").Value = -20320.5386

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")

# Row 5
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 329.83334
$arr[0,1] = 297.25
$arr[0,2] = 395
$arr[0,3] = 297.25
$arr[0,4] = 395
$arr[0,5] = -185.25
$arr[0,6] = -619
$ws.Range("H5:N5").Value = $arr

# Row 32
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 5478.71
$arr[0,1] = 4298.6665
$arr[0,2] = 21156.428
$arr[0,3] = 4298.6665
$arr[0,4] = 21156.428
$arr[0,5] = -4011.6665
$arr[0,6] = -21730.428
$ws.Range("H32:N32").Value = $arr

# Row 61
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2447
$arr[0,1] = 1435.8182
$arr[0,2] = 6154.6665
$arr[0,3] = 1435.8182
$arr[0,4] = 6154.6665
$arr[0,5] = -1223.8182
$arr[0,6] = -6578.6665
$ws.Range("H61:N61").Value = $arr

# Row 74
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2877.8076
$arr[0,1] = 2013.4348
$arr[0,2] = 9504.666999999999
$arr[0,3] = 2013.4348
$arr[0,4] = 9504.666999999999
$arr[0,5] = -1139.4348
$arr[0,6] = -11252.667
$ws.Range("H74:N74").Value = $arr

# Row 77
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2877.8076
$arr[0,1] = 2013.4348
$arr[0,2] = 9504.666999999999
$arr[0,3] = 10067.174
$arr[0,4] = 47523.335
$arr[0,5] = -5699.173999999999
$arr[0,6] = -56259.335
$ws.Range("H77:N77").Value = $arr

# Row 80
$ws.Range("H80").Value = 45422
$ws.Range("J80").Value = 48527.5
$ws.Range("L80").Value = 48527.5
$ws.Range("N80").Value = -50523.5

# Row 83
$ws.Range("H83").Value = 45422
$ws.Range("J83").Value = 48527.5
$ws.Range("L83").Value = 145582.5
$ws.Range("N83").Value = -155566.5

# Row 132
$ws.Range("H132").Value = 3666.8
$ws.Range("I132").Value = 3552.25
$ws.Range("K132").Value = 10656.75
$ws.Range("M132").Value = -8126.75

# Row 136
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2447
$arr[0,1] = 1435.8182
$arr[0,2] = 6154.6665
$arr[0,3] = 4307.4546
$arr[0,4] = 18463.9995
$arr[0,5] = -1757.4546
$arr[0,6] = -23563.9995
$ws.Range("H136:N136").Value = $arr

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")

# Row 4
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 329.83334
$arr[0,1] = 297.25
$arr[0,2] = 395
$arr[0,3] = 297.25
$arr[0,4] = 395
$arr[0,5] = -182.25
$arr[0,6] = -625
$ws.Range("H4:N4").Value = $arr

# Row 134
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 3626.75
$arr[0,1] = 3419
$arr[0,2] = 4250
$arr[0,3] = 10257
$arr[0,4] = 12750
$arr[0,5] = -7722
$arr[0,6] = -17820
$ws.Range("H134:N134").Value = $arr

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")

# Row 7
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 13473.134
$arr[0,1] = 12645.875
$arr[0,2] = 14418.571
$arr[0,3] = 12645.875
$arr[0,4] = 14418.571
$arr[0,5] = -12532.875
$arr[0,6] = -14644.571
$ws.Range("H7:N7").Value = $arr

# Row 31
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 4841.4
$arr[0,1] = 2809.7896
$arr[0,2] = 6086.5806
$arr[0,3] = 2809.7896
$arr[0,4] = 6086.5806
$arr[0,5] = -2514.7896
$arr[0,6] = -6676.5806
$ws.Range("H31:N31").Value = $arr

# Row 34
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 4841.4
$arr[0,1] = 2809.7896
$arr[0,2] = 6086.5806
$arr[0,3] = 2809.7896
$arr[0,4] = 6086.5806
$arr[0,5] = -2607.7896
$arr[0,6] = -6490.5806
$ws.Range("H34:N34").Value = $arr

# Row 99
$ws.Range("H99").Value = 2085873.5
$ws.Range("I99").Value = 2068.3333
$ws.Range("K99").Value = 2068.3333
$ws.Range("M99").Value = -570.3332999999998

# Row 105
$ws.Range("H105").Value = 94995.5
$ws.Range("I105").Value = 224231.6
$ws.Range("K105").Value = 224231.6
$ws.Range("M105").Value = -222484.6

# Row 126
$ws.Range("H126").Value = 2085873.5
$ws.Range("I126").Value = 2068.3333
$ws.Range("K126").Value = 6204.999899999999
$ws.Range("M126").Value = -3734.999899999999

# Row 132
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2802.4
$arr[0,1] = 3021
$arr[0,2] = 2708.7144
$arr[0,3] = 9063
$arr[0,4] = 8126.1432
$arr[0,5] = -6533
$arr[0,6] = -13186.1432
$ws.Range("H132:N132").Value = $arr

# Row 134
$ws.Range("H134").Value = 1607.5294
$ws.Range("I134").Value = 1339.8928
$ws.Range("K134").Value = 4019.6784
$ws.Range("M134").Value = -1484.6784

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")

# Row 38
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 40.727272
$arr[0,1] = 34.714287
$arr[0,2] = 43.533333
$arr[0,3] = 104.142861
$arr[0,4] = 130.599999
$arr[0,5] = 242.857139
$arr[0,6] = -824.599999
$ws.Range("H38:N38").Value = $arr

# Row 86
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 3191.5833
$arr[0,1] = 1542.7142
$arr[0,2] = 3870.5293
$arr[0,3] = 4628.142599999999
$arr[0,4] = 11611.5879
$arr[0,5] = -3442.142599999999
$arr[0,6] = -13983.5879
$ws.Range("H86:N86").Value = $arr

# Row 89
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 3191.5833
$arr[0,1] = 1542.7142
$arr[0,2] = 3870.5293
$arr[0,3] = 13884.4278
$arr[0,4] = 34834.7637
$arr[0,5] = -7956.427799999999
$arr[0,6] = -46690.7637
$ws.Range("H89:N89").Value = $arr

# Row 131
$ws.Range("H131").Value = 29918.416
$ws.Range("J131").Value = 2450.1538
$ws.Range("L131").Value = 7350.4614
$ws.Range("N131").Value = -17430.4614

# Row 140
$ws.Range("H140").Value = 2418.7693
$ws.Range("I140").Value = 1749.2858
$ws.Range("K140").Value = 5247.857400000001
$ws.Range("M140").Value = -67.85740000000078

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")

# Row 132
$ws.Range("H132").Value = 2360.9092
$ws.Range("I132").Value = 1969
$ws.Range("K132").Value = 5907
$ws.Range("M132").Value = -3377

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")

# Row 7
$ws.Range("H7").Value = 13443.435
$ws.Range("I7").Value = 12510.6
$ws.Range("K7").Value = 12510.6
$ws.Range("M7").Value = -12398.6

# Row 31
$ws.Range("H31").Value = 4658.143
$ws.Range("J31").Value = 6402.8
$ws.Range("L31").Value = 6402.8
$ws.Range("N31").Value = -6898.8

# Row 34
$ws.Range("H34").Value = 27499.5
$ws.Range("J34").Value = 27499.5
$ws.Range("L34").Value = 27499.5
$ws.Range("N34").Value = -27843.5

# Row 43
$ws.Range("H43").Value = 16000
$ws.Range("J43").Value = 16000
$ws.Range("L43").Value = 16000
$ws.Range("N43").Value = -16386

# Row 46
$ws.Range("H46").Value = 5841.52
$ws.Range("I46").Value = 10719.909
$ws.Range("K46").Value = 10719.909
$ws.Range("M46").Value = -10531.909

# Row 126
$ws.Range("H126").Value = 13443.435
$ws.Range("I126").Value = 12510.6
$ws.Range("K126").Value = 37531.8
$ws.Range("M126").Value = -35061.8

# Row 136
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 9939.4
$arr[0,1] = 10715.1
$arr[0,2] = 8388
$arr[0,3] = 32145.3
$arr[0,4] = 25164
$arr[0,5] = -29595.3
$arr[0,6] = -30264
$ws.Range("H136:N136").Value = $arr

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")

# Row 132
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 10876565
$arr[0,1] = 11999.5
$arr[0,2] = 21741130
$arr[0,3] = 35998.5
$arr[0,4] = 65223390
$arr[0,5] = -33468.5
$arr[0,6] = -65228450
$ws.Range("H132:N132").Value = $arr

# Row 136
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2531.3171
$arr[0,1] = 1803
$arr[0,2] = 6068.857
$arr[0,3] = 5409
$arr[0,4] = 18206.571
$arr[0,5] = -2859
$arr[0,6] = -23306.571
$ws.Range("H136:N136").Value = $arr
